$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the old row 85 (pushing old rows 85-115 down to 88-118).
# Excel's Insert() carries the formatting of the row above (row 84) onto the new rows,
# which matches styles s=54 (col A), s=52 (cols B-E), s=5 (F), s=6 (G/H/N-S), s=7 (I),
# s=8 (K), s=9 (L), s=10 (M) seen on row 84.
$ws.Rows("85:87").Insert()

# New accession / sequence IDs (column A).
$ws.Range("A85").Value = "MH716818"
$ws.Range("A86").Value = "MG599986"
$ws.Range("A87").Value = "MF776369"

# New virus full names (column B).
$ws.Range("B85").Value = "Eastern red scorpionfish flavivirus"
$ws.Range("B86").Value = "Wenzhou shark flavivirus"

# New sampled host scientific names (column G).
$ws.Range("G86").Value = "Scoliodon macrorhynchos"
$ws.Range("G85").Value = "Scorpaena jacksoniensis"

$ws.Range("B87").Value = "Cyclopterus lumpus virus"
$ws.Range("G87").Value = "Cyclopterus lumpus"

# Re-point the frozen pane from a column freeze (freeze col A) to a row freeze
# (freeze row 1), and scroll/select down near the newly-added rows.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true

$ws.Range("F87").Select()
